# Update "想去人数" (F column) figures on sheets "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheet1Updates = @{
    2  = 278
    3  = 566
    5  = 274
    6  = 1078
    7  = 1411
    9  = 102
    10 = 740
    11 = 67
    12 = 135
    13 = 124
    14 = 420
    15 = 1320
    16 = 102
    17 = 89
    18 = 269
    19 = 5214
    20 = 641
    21 = 30
    22 = 202
    24 = 5688
    25 = 56
    26 = 118
    29 = 14247
    30 = 1420
    31 = 197
    32 = 93
    34 = 432
    35 = 586
    36 = 4178
}

$sheet4Updates = @{
    2  = 278
    3  = 566
    5  = 274
    6  = 1078
    7  = 1411
    9  = 102
    10 = 740
    11 = 67
    12 = 135
    13 = 124
    14 = 420
    15 = 1320
    16 = 102
    17 = 89
    18 = 269
    20 = 5214
    21 = 641
    23 = 30
    24 = 202
    25 = 13
    27 = 5688
    28 = 56
    29 = 118
    32 = 14247
    33 = 1420
    34 = 197
    35 = 93
    37 = 432
    38 = 586
    39 = 4178
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
